# Update newsbot state: drop the (empty) F44/G44 placeholder cells and
# append the newly scraped article as row 45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 loses its trailing empty F/G cells.
$ws.Range("F44").ClearContents()
$ws.Range("G44").ClearContents()

# New row 45 with the latest news item.
$ws.Range("A45").Value = "05/01/2026 06:47:41"
$ws.Range("B45").Value = "05/01 06:33"
$ws.Range("C45").Value = "Metrópoles"
$ws.Range("D45").Value = "Lula volta do recesso com foco na Venezuela, mas tem outras pendências"
$ws.Range("E45").Value = "https://www.metropoles.com/brasil/lula-volta-do-recesso-com-foco-na-venezuela-mas-tem-outras-pendencias"
$ws.Range("F45").Value = "senado"
$ws.Range("G45").Value = "Presidente precisa conversar com o presidente do Senado, Davi Alcolumbre, e decidir mudanças no ministério de Lewandowski"
